# Updated symbol list on Sat Dec 24 10:58:25 UTC 2022 with GitHub Actions
#
# This refreshes the crypto price/volume snapshot on the sheet. Every cell
# in this sheet is stored as text (prices like "244.98" are inline strings,
# not numbers), so values are written as text and the cached number format
# is cleared afterwards so Excel doesn't silently reinterpret/restyle the
# numeric-looking strings as real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.ClearFormats()
}

# --- Price (column D) corrections across the table ---
Set-Text "D2"  "244.82"
Set-Text "D3"  "22.00"
Set-Text "D4"  "5.400"
Set-Text "D5"  "0.06003"
Set-Text "D7"  "0.8128"
Set-Text "D8"  "0.9578"
Set-Text "D10" "0.07358"
Set-Text "D12" "0.03054"
Set-Text "D14" "4.001"
Set-Text "D15" "0.001590"
Set-Text "D16" "0.04813"
Set-Text "D17" "0.0005873"
Set-Text "E17" "16OneONEWorstin24h"
Set-Text "D18" "0.006234"
Set-Text "D19" "0.005044"
Set-Text "D20" "0.0009904"
Set-Text "D22" "3.698"
Set-Text "D23" "6.423"
Set-Text "D26" "0.1340"
Set-Text "D40" "0.03998"

# --- Rows 41-43: ranking reshuffled (KickToken moved up, BKEXToken and CEJI
#     shifted down one place each) ---
Set-Text "B41" "KickToken"
Set-Text "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-Text "D41" "0.006543"
Set-Text "E41" "40KickTokenKICK"

Set-Text "B42" "BKEXToken"
Set-Text "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-Text "D42" "0.1072"
Set-Text "E42" "41BKEXTokenBKK"

Set-Text "B43" "CEJI"
Set-Text "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-Text "D43" "0.002901"
Set-Text "E43" "42CEJICEJI"

Set-Text "D44" "0.005842"
Set-Text "D45" "0.00005256"
Set-Text "D47" "0.9694"
Set-Text "D48" "0.02129"
Set-Text "E48" "47BOLOBOLO"
